# Preparing for second keyboard
# Applies the changes described by the target diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear the old Windows/Premiere hotkey notes out of rows 3-7 (columns
#    E and G) - they are being relocated into the new "Second Keyboard"
#    block further down the sheet.
# ---------------------------------------------------------------------
$ws.Range("E3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("G7").ClearContents()

# ---------------------------------------------------------------------
# 2. Add a "reset" label above the livesplit undo/skip/split list (I49).
# ---------------------------------------------------------------------
$ws.Range("I49").Value = "reset"
$ws.Range("I49").Style = "Good"

# ---------------------------------------------------------------------
# 3. New column width for column C (so the "Second Keyboard" label
#    fits).
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.42578125

# ---------------------------------------------------------------------
# 4. Build the new "Second Keyboard" block (rows 89-142), mirroring the
#    existing Numpad/Other/Other-hotkeys sections one column to the
#    right.
# ---------------------------------------------------------------------

# Section header + numpad label
$ws.Range("C89").Value = "Second Keyboard"
$ws.Range("D89").Value = "Numpad"
$ws.Range("D89").Style = "Normal"

# Numpad digits 1-9 in column E, rows 89-97 (left aligned, like the
# original numpad block)
for ($i = 0; $i -le 8; $i++) {
    $r = 89 + $i
    $cell = $ws.Range("E$r")
    $cell.Value = $i + 1
    $cell.HorizontalAlignment = -4131   # xlLeft
}

# Mouse/drag related notes moved from the old rows 3-7 (column G),
# realigned with rows 93-97
$ws.Range("G93").Value = "adjust rotation"
$ws.Range("G93").Style = "Good"
$ws.Range("G94").Value = "drag via preview window"
$ws.Range("G94").Style = "Good"
$ws.Range("G95").Value = "adjust scale"
$ws.Range("G95").Style = "Good"
$ws.Range("G96").Value = "move x axis"
$ws.Range("G96").Style = "Good"
$ws.Range("G97").Value = "move y axis"
$ws.Range("G97").Style = "Good"

# Remaining numpad keys
$ws.Range("E98").Value = "/"
$ws.Range("D98").Style = "Normal"

$ws.Range("E99").Value = "-"
$ws.Range("G99").Value = "explorer"
$ws.Range("G99").Style = "Good"

$ws.Range("E100").Value = "+"
$ws.Range("G100").Value = "premiere"
$ws.Range("G100").Style = "Good"

$ws.Range("E101").Value = "NumLock"

$ws.Range("E102").Value = "backspace"
$ws.Range("G102").Value = "firefox"
$ws.Range("G102").Style = "Good"

# blank spacer row
$ws.Range("D103").Style = "Normal"

# blank row before the "Other" block
$ws.Range("D104").Style = "Normal"
$ws.Range("E104").Style = "Normal"

# "Other" block (Rwin / AppsKey / Ctrl AppsKey), shifted one column
# right compared to the original C70:E72 block
$ws.Range("D105").Value = "Other"
$ws.Range("E105").Value = "Rwin"
$ws.Range("D106").Style = "Normal"
$ws.Range("E106").Value = "AppsKey"
$ws.Range("D107").Style = "Normal"
$ws.Range("E107").Value = "Ctrl AppsKey"

# blank rows between "Other" and "Other hotkeys"
$ws.Range("D108").Style = "Normal"
$ws.Range("E108").Style = "Normal"
$ws.Range("D109").Style = "Normal"
$ws.Range("E109").Style = "Normal"
$ws.Range("D110").Style = "Normal"
$ws.Range("E110").Style = "Normal"
$ws.Range("D111").Style = "Normal"
$ws.Range("E111").Style = "Normal"

# "Other hotkeys" label
$ws.Range("D112").Style = "Normal"
$ws.Range("E112").Value = "Other hotkeys"

# trailing blank rows 113-142 (kept styled/blank, matching the sheet's
# used range extending down to row 142)
for ($r = 113; $r -le 142; $r++) {
    $ws.Range("D$r").Style = "Normal"
    $ws.Range("E$r").Style = "Normal"
}

# ---------------------------------------------------------------------
# 5. Update the view so it lands where the author left it (scrolled to
#    the new block, with H77 selected).
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("H77").Select()
